# Weekly update: a new "Jengibre" (ginger) price-report row for
# Terminal La Palmera de La Serena is inserted at row 15, pushing the
# existing rows 15-35 down to 16-36 (dimension grows from R35 to R36).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 15; rows 15..35 shift to 16..36.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with this week's record.
$ws.Cells.Item(15, 1).Value = 8
$ws.Cells.Item(15, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(15, 3).Value = "Coquimbo"
$ws.Cells.Item(15, 4).Value = 44658
$ws.Cells.Item(15, 5).Value = 4
$ws.Cells.Item(15, 6).Value = 100114007
$ws.Cells.Item(15, 7).Value = "Jengibre"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 400
$ws.Cells.Item(15, 11).Value = 15000
$ws.Cells.Item(15, 12).Value = 16000
$ws.Cells.Item(15, 13).Value = 15500
$ws.Cells.Item(15, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(15, 15).Value = "Perú"
$ws.Cells.Item(15, 16).Value = 1192
$ws.Cells.Item(15, 17).Value = 13
$ws.Cells.Item(15, 18).Value = "Hortaliza"
